$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New student data replacing the old two rows (Id, Name, Age, Year)
$students = @(
    @(1, "Swen",    19, 3),
    @(2, "Anthony", 19, 3),
    @(3, "Hans",    25, 4),
    @(4, "Frank",   21, 4),
    @(5, "Jan",     23, 2)
)

$row = 11
foreach ($s in $students) {
    $ws.Cells.Item($row, 10).Value = $s[0]   # J - Id
    $ws.Cells.Item($row, 11).Value = $s[1]   # K - Name
    $ws.Cells.Item($row, 12).Value = $s[2]   # L - Age
    $ws.Cells.Item($row, 13).Value = $s[3]   # M - Year
    $row++
}
